$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date text (18/05/2020 -> 21/05/2020)
#    on the slide master and every slide layout footer placeholder.
# ---------------------------------------------------------------------------
function Update-DateField($shapes, [string]$oldText, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

Update-DateField $p.SlideMaster.Shapes "18/05/2020" "21/05/2020"

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes "18/05/2020" "21/05/2020"
}

# ---------------------------------------------------------------------------
# 2) Add a new slide 2 ("Titre seul" layout) with the title "Cleaning" and a
#    7x8 results table.
# ---------------------------------------------------------------------------
$titleOnlyLayout = $p.SlideMaster.CustomLayouts.Item(6)
$s2 = $p.Slides.Add(2, $titleOnlyLayout)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Cleaning"

# EMU -> point conversion (PowerPoint Shapes.AddTable takes points).
$emuPerPt = 12700.0
$offX = 467544 / $emuPerPt
$offY = 1412776 / $emuPerPt
$extCX = 8280920 / $emuPerPt
$extCY = 4300008 / $emuPerPt

$tblShape = $s2.Shapes.AddTable(7, 8, $offX, $offY, $extCX, $extCY)
$tblShape.Name = "Tableau 2"
$tbl = $tblShape.Table
$tbl.ApplyStyle("{5940675A-B579-460E-94D1-54222C63F5DA}")

$colWidths = @(1035115, 1035115, 1035115, 1035115, 1044116, 1026114, 1035115, 1035115)
for ($c = 1; $c -le 8; $c++) {
    $tbl.Columns.Item($c).Width = $colWidths[$c - 1] / $emuPerPt
}

$rowHeights = @(792088, 576064, 576064, 588948, 588948, 588948, 588948)
for ($r = 1; $r -le 7; $r++) {
    $tbl.Rows.Item($r).Height = $rowHeights[$r - 1] / $emuPerPt
}

# Header row content
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Text = "Type"
$tbl.Cell(1, 4).Shape.TextFrame.TextRange.Text = "Nb deleted"
$tbl.Cell(1, 5).Shape.TextFrame.TextRange.Text = "Nb fixed"
$tbl.Cell(1, 7).Shape.TextFrame.TextRange.Text = "Nb col"
$tbl.Cell(1, 8).Shape.TextFrame.TextRange.Text = "nbcol"

# Row index column + lone "columns" note
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "0"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "1"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "columns"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "2"
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "3"
$tbl.Cell(6, 1).Shape.TextFrame.TextRange.Text = "4"
$tbl.Cell(7, 1).Shape.TextFrame.TextRange.Text = "5"

Write-Output "done"
